# Apply cryptos list update (price/volume refresh + two row swaps)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '22.428.77'
$ws.Range("E2").Value = '  +8.81%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.603.50'
$ws.Range("E3").Value = '  +8.38%  '

$ws.Range("E4").Value = '  -0.72%  '

$ws.Range("B5").Value = 'USDC'
$ws.Range("C5").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.9915'
$ws.Range("E5").Value = '  +3.01%  '

$ws.Range("B6").Value = 'BNB'
$ws.Range("C6").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '303.11'
$ws.Range("E6").Value = '  +7.94%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3687'
$ws.Range("E7").Value = '  +0.56%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3395'
$ws.Range("E8").Value = '  +10.08%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '42.36'
$ws.Range("E9").Value = '  +5.78%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.138'
$ws.Range("E10").Value = '  +7.06%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07060'
$ws.Range("E11").Value = '  +5.82%  '

$ws.Range("E12").Value = '  -0.46%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.70'
$ws.Range("E13").Value = '  +8.73%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.924'
$ws.Range("E14").Value = '  +7.08%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.640'
$ws.Range("E15").Value = '  +6.74%  '

$ws.Range("E16").Value = '  +5.04%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.602.77'
$ws.Range("E17").Value = '  +8.45%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.9916'
$ws.Range("E18").Value = '  +2.78%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06859'
$ws.Range("E19").Value = '  +14.80%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '78.01'
$ws.Range("E20").Value = '  +11.26%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.046'
$ws.Range("E21").Value = '  +9.83%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '16.11'
$ws.Range("E22").Value = '  +11.17%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.83'
$ws.Range("E23").Value = '  +6.87%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '22.460.91'
$ws.Range("E24").Value = '  +8.78%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.382'
$ws.Range("E25").Value = '  +5.46%  '

$ws.Range("E26").Value = '  +19.56%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '151.59'
$ws.Range("E27").Value = '  +6.16%  '

$ws.Range("E28").Value = '  +13.36%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.783.33'
$ws.Range("E29").Value = '  +8.79%  '

$ws.Range("E30").Value = '  +6.00%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.139'
$ws.Range("E31").Value = '  +4.06%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.072'
$ws.Range("E32").Value = '  +20.08%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.9523'
$ws.Range("E33").Value = '  +16.07%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08266'
$ws.Range("E34").Value = '  +3.26%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.639'
$ws.Range("E35").Value = '  +6.78%  '

$ws.Range("B36").Value = 'Aptos'
$ws.Range("C36").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '12.02'
$ws.Range("E36").Value = '  +15.22%  '

$ws.Range("B37").Value = 'InternetComputer(DFINITY)'
$ws.Range("C37").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.266'
$ws.Range("E37").Value = '  +10.96%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.619'
$ws.Range("E38").Value = '  +14.30%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.259'
$ws.Range("E39").Value = '  +3.21%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.06099'
$ws.Range("E40").Value = '  +4.84%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.02222'
$ws.Range("E41").Value = '  +8.31%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.2022'
$ws.Range("E42").Value = '  +7.44%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9918'
$ws.Range("E43").Value = '  +2.88%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5922'
$ws.Range("E44").Value = '  +11.28%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.819'
$ws.Range("E45").Value = '  +7.66%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '13.20'
$ws.Range("E46").Value = '  +6.73%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5706'
$ws.Range("E47").Value = '  +9.40%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '127.66'
$ws.Range("E48").Value = '  +7.60%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.983'
$ws.Range("E49").Value = '  +8.17%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06801'
$ws.Range("E50").Value = '  +4.55%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '73.79'
$ws.Range("E51").Value = '  +8.76%  '
